# Update the result score for HoangNC from "18/34" to "24/34".
#
# The target text lives in a table cell as a single run:
#   "18/34"  (rFonts=MS Gothic/minorHAnsi, bold, color FF0000)
#
# After the edit it must become two runs with identical formatting:
#   "24"  +  "/34"
#
# We locate "18/34" with Find, overwrite just the "18" portion with "24",
# then round-trip a Font property on the newly written "24" text so the
# engine materialises it as its own run instead of silently re-merging it
# with the untouched "/34" text that follows (both runs end up with the
# same rPr, exactly like the target diff).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("18/34", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '18/34' in the document"
}

$matchStart = $find.Parent.Start
$matchEnd = $find.Parent.End
if ($d.Range($matchStart, $matchEnd).Text -ne "18/34") {
    throw "Unexpected match text for '18/34'"
}

# "18" is the first two characters of the "18/34" run; "/34" is the rest
# and must stay untouched.
$oldScorePart = $d.Range($matchStart, $matchStart + 2)
$oldScorePart.Text = "24"

# Force a run boundary right after the new "24" by round-tripping the
# font color (read it, change it, set it back) on exactly that range.
$newScorePart = $d.Range($matchStart, $matchStart + 2)
$originalColor = $newScorePart.Font.Color
$newScorePart.Font.Color = $originalColor + 1
$newScorePart.Font.Color = $originalColor

Write-Output "Updated score: $($d.Range($matchStart, $matchStart + 5).Text)"
